$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete now-unused trailing rows (41-46) so dimension shrinks to A1:C40
$ws.Rows("41:46").Delete()

# Update rows 2-40 with new "#", "Song", "Registration #" values
$ws.Cells.Item(2, 1).Value = 1
$ws.Cells.Item(2, 2).Value = '13463-20th Century Boy'
$ws.Cells.Item(2, 3).Value = 'PA0000976753'
$ws.Cells.Item(3, 1).Value = 2
$ws.Cells.Item(3, 2).Value = '54307-WHO ARE YOU'
$ws.Cells.Item(3, 3).Value = 'PA0001864711'
$ws.Cells.Item(4, 1).Value = 3
$ws.Cells.Item(4, 2).Value = '56418-BABA O''RILEY'
$ws.Cells.Item(4, 3).Value = 'PA0000720226'
$ws.Cells.Item(5, 1).Value = 4
$ws.Cells.Item(5, 2).Value = '191848-MEANT TO BE'
$ws.Cells.Item(5, 3).Value = 'PA0001751784'
$ws.Cells.Item(6, 1).Value = 8
$ws.Cells.Item(6, 2).Value = '191897-SOUTHBOUND'
$ws.Cells.Item(6, 3).Value = 'PA0002267459'
$ws.Cells.Item(7, 1).Value = 9
$ws.Cells.Item(7, 2).Value = '202993-BEER NEVER BROKE MY HEART'
$ws.Cells.Item(7, 3).Value = 'PA0002192567'
$ws.Cells.Item(8, 1).Value = 11
$ws.Cells.Item(8, 2).Value = '24776-CAN''T YOU SEE U.S. ONLY AS OF)'
$ws.Cells.Item(8, 3).Value = 'SR0000863063'
$ws.Cells.Item(9, 1).Value = 12
$ws.Cells.Item(9, 2).Value = '19819-20th Century Boy - Master'
$ws.Cells.Item(9, 3).Value = 'SR0000233398'
$ws.Cells.Item(10, 1).Value = 14
$ws.Cells.Item(10, 2).Value = '11872-Heaven Is A Place On Earth'
$ws.Cells.Item(10, 3).Value = 'PA0000877899'
$ws.Cells.Item(11, 1).Value = 15
$ws.Cells.Item(11, 2).Value = '56421-BEHIND BLUE EYES'
$ws.Cells.Item(11, 3).Value = 'PA0000347339'
$ws.Cells.Item(12, 1).Value = 16
$ws.Cells.Item(12, 2).Value = '56531-PINBALL WIZARD'
$ws.Cells.Item(12, 3).Value = 'PA0000750495'
$ws.Cells.Item(13, 1).Value = 17
$ws.Cells.Item(13, 2).Value = '13262-Children Of The Revolution'
$ws.Cells.Item(13, 3).Value = 'PA0001015278'
$ws.Cells.Item(14, 1).Value = 19
$ws.Cells.Item(14, 2).Value = '191840-LOVE WINS'
$ws.Cells.Item(14, 3).Value = 'PA0002031121'
$ws.Cells.Item(15, 1).Value = 21
$ws.Cells.Item(15, 2).Value = '1943-GONNA MAKE YOU SWEAT'
$ws.Cells.Item(15, 3).Value = 'PA0000520237'
$ws.Cells.Item(16, 1).Value = 22
$ws.Cells.Item(16, 2).Value = '100056-TEACH YOUR CHILDREN'
$ws.Cells.Item(16, 3).Value = 'V3615D973'
$ws.Cells.Item(17, 1).Value = 23
$ws.Cells.Item(17, 2).Value = '12008-You Get What You Give'
$ws.Cells.Item(17, 3).Value = 'V3615D973'
$ws.Cells.Item(18, 1).Value = 24
$ws.Cells.Item(18, 2).Value = '191644-DRINKING ALONE'
$ws.Cells.Item(18, 3).Value = 'PA0002233922'
$ws.Cells.Item(19, 1).Value = 25
$ws.Cells.Item(19, 2).Value = '182638-BODY LIKE A BACK ROAD'
$ws.Cells.Item(19, 3).Value = 'PA0002067268'
$ws.Cells.Item(20, 1).Value = 27
$ws.Cells.Item(20, 2).Value = '56635-MY GENERATION'
$ws.Cells.Item(20, 3).Value = 'PA0000722088'
$ws.Cells.Item(21, 1).Value = 28
$ws.Cells.Item(21, 2).Value = '165764-THIS IS IT'
$ws.Cells.Item(21, 3).Value = 'PA0001618358'
$ws.Cells.Item(22, 1).Value = 30
$ws.Cells.Item(22, 2).Value = '56469-THE SEEKER'
$ws.Cells.Item(22, 3).Value = 'PA0002029207'
$ws.Cells.Item(23, 1).Value = 31
$ws.Cells.Item(23, 2).Value = '190499-THE WAY I AM'
$ws.Cells.Item(23, 3).Value = 'PA0001601324'
$ws.Cells.Item(24, 1).Value = 32
$ws.Cells.Item(24, 2).Value = '56584-JOIN TOGETHER'
$ws.Cells.Item(24, 3).Value = 'PA0001218253'
$ws.Cells.Item(25, 1).Value = 33
$ws.Cells.Item(25, 2).Value = '149616-CHA CHA SLIDE'
$ws.Cells.Item(25, 3).Value = 'PA0001204600'
$ws.Cells.Item(26, 1).Value = 35
$ws.Cells.Item(26, 2).Value = '12188-DANCE HALL DAYS'
$ws.Cells.Item(26, 3).Value = 'PA0001074195'
$ws.Cells.Item(27, 1).Value = 37
$ws.Cells.Item(27, 2).Value = 'NL50499-TOMMY - STAGE INCOME'
$ws.Cells.Item(27, 3).Value = 'PA0000539941'
$ws.Cells.Item(28, 1).Value = 38
$ws.Cells.Item(28, 2).Value = '56638-I''M A BOY'
$ws.Cells.Item(28, 3).Value = 'PA0001218257'
$ws.Cells.Item(29, 1).Value = 39
$ws.Cells.Item(29, 2).Value = '12203-EVERYBODY HAVE FUN TONIGHT'
$ws.Cells.Item(29, 3).Value = 'PA0000324723'
$ws.Cells.Item(30, 1).Value = 40
$ws.Cells.Item(30, 2).Value = '480-CHRISTMAS WRAPPING'
$ws.Cells.Item(30, 3).Value = ' '
$ws.Cells.Item(31, 1).Value = 41
$ws.Cells.Item(31, 2).Value = '193999-HERE AND NOW'
$ws.Cells.Item(31, 3).Value = 'PA0000394195'
$ws.Cells.Item(32, 1).Value = 42
$ws.Cells.Item(32, 2).Value = '10841-FAMILY AFFAIR BOTH SHARES)'
$ws.Cells.Item(32, 3).Value = ' '
$ws.Cells.Item(33, 1).Value = 43
$ws.Cells.Item(33, 2).Value = '92403-BAREFOOT BLUE JEAN NIGHT BMI'
$ws.Cells.Item(33, 3).Value = 'SR0000697851'
$ws.Cells.Item(34, 1).Value = 44
$ws.Cells.Item(34, 2).Value = '182982-SUNRISE SUNBURN SUNSET'
$ws.Cells.Item(34, 3).Value = 'PA0002140607'
$ws.Cells.Item(35, 1).Value = 45
$ws.Cells.Item(35, 2).Value = '187657-GOOD VIBES'
$ws.Cells.Item(35, 3).Value = 'PA0002228161'
$ws.Cells.Item(36, 1).Value = 46
$ws.Cells.Item(36, 2).Value = '187820-RIDIN'' ROADS'
$ws.Cells.Item(36, 3).Value = 'PA0002197474'
$ws.Cells.Item(37, 1).Value = 47
$ws.Cells.Item(37, 2).Value = '203068-I HOPE YOU''RE HAPPY NOW'
$ws.Cells.Item(37, 3).Value = 'PA0002297072'
$ws.Cells.Item(38, 1).Value = 48
$ws.Cells.Item(38, 2).Value = '1522-LET''S TWIST AGAIN'
$ws.Cells.Item(38, 3).Value = 'PA0000196738'
$ws.Cells.Item(39, 1).Value = 49
$ws.Cells.Item(39, 2).Value = '20839-I Know You Want Me'
$ws.Cells.Item(39, 3).Value = 'TX0002685123'
$ws.Cells.Item(40, 1).Value = 50
$ws.Cells.Item(40, 2).Value = '54313-EMINENCE FRONT'
$ws.Cells.Item(40, 3).Value = 'PA0000152030'
